# "Colocando header nos gráficos" - add header labels to the data tables
# that feed the charts, fix accentuation typos, and update a couple of
# values on the "Custo Total" sheet. Also drop the now-unused "Teto" row
# on the "Emissoes Totais" sheet.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell {
    param($ws, [string]$addr, [string]$text)
    # Give the new header cell the same look (bold / border / centered)
    # as the existing header row by copying format from B1, which
    # already carries that style.
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value = $text
}

function Clear-CellStyle {
    param($ws, [string]$addr)
    $ws.Range($addr).ClearFormats()
}

# ---------------------------------------------------------------
# Sheets 1-4 share the same row layout (Fonte/Tecnologia column).
# ---------------------------------------------------------------
$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header for column A
    Set-HeaderCell $ws "A1" "Fonte/Tecnologia"

    # Row labels lose their bold/border style and a few get accents fixed
    Clear-CellStyle $ws "A2"
    $ws.Range("A2").Value = "Hidro"

    Clear-CellStyle $ws "A3"
    $ws.Range("A3").Value = "Gás Natural"

    Clear-CellStyle $ws "A4"
    $ws.Range("A4").Value = "Carvão"

    Clear-CellStyle $ws "A5"
    $ws.Range("A5").Value = "Nuclear"

    Clear-CellStyle $ws "A6"
    $ws.Range("A6").Value = "Óleos Comb"

    Clear-CellStyle $ws "A7"
    $ws.Range("A7").Value = "Biomassa"

    Clear-CellStyle $ws "A8"
    $ws.Range("A8").Value = "Eólica"

    Clear-CellStyle $ws "A9"
    $ws.Range("A9").Value = "Solar"

    Clear-CellStyle $ws "A10"
    $ws.Range("A10").Value = "Outros"

    Clear-CellStyle $ws "A11"
    $ws.Range("A11").Value = "Pot. Compl."

    Clear-CellStyle $ws "A12"
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)" - add header, fix labels,
# and remove the unused "Teto" row (row 4).
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Rows.Item(4).Delete() | Out-Null

Set-HeaderCell $ws5 "A1" "Período"

Clear-CellStyle $ws5 "A2"
$ws5.Range("A2").Value = "P.Médio"

Clear-CellStyle $ws5 "A3"
$ws5.Range("A3").Value = "P.Crítico"

# ---------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)" - add header, fix labels,
# update the 2015 cost figures.
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

Set-HeaderCell $ws6 "A1" "Tipo Expansão"

# B1 needs to become the text "2015" (same header text/style used on the
# other sheets) rather than a number, so copy the whole cell (value +
# format) from a sheet that already has a "2015" text header.
$ws1 = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
$ws1.Range("B1").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4104) | Out-Null   # xlPasteAll

Clear-CellStyle $ws6 "A2"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 592

Clear-CellStyle $ws6 "A3"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
